$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 1-63 (header + 62 data rows)
$data = @(
    ,@('Contract Name', 'Dependencies', 'Dependent Contracts')
    ,@('IUniswapV3FlashCallback', 0, 0)
    ,@('IUniswapV3MintCallback', 0, 0)
    ,@('IUniswapV3SwapCallback', 0, 0)
    ,@('IERC20Minimal', 0, 0)
    ,@('IUniswapV3Factory', 0, 0)
    ,@('IUniswapV3Pool', 6, 0)
    ,@('IUniswapV3PoolDeployer', 0, 0)
    ,@('IUniswapV3PoolActions', 0, 0)
    ,@('IUniswapV3PoolDerivedState', 0, 0)
    ,@('IUniswapV3PoolEvents', 0, 0)
    ,@('IUniswapV3PoolImmutables', 0, 0)
    ,@('IUniswapV3PoolOwnerActions', 0, 0)
    ,@('IUniswapV3PoolState', 0, 0)
    ,@('BitMath', 0, 0)
    ,@('FixedPoint128', 0, 0)
    ,@('FixedPoint96', 0, 0)
    ,@('FullMath', 0, 0)
    ,@('LiquidityMath', 0, 0)
    ,@('LowGasSafeMath', 0, 0)
    ,@('Oracle', 0, 0)
    ,@('Position', 3, 0)
    ,@('SafeCast', 0, 0)
    ,@('SqrtPriceMath', 5, 0)
    ,@('SwapMath', 2, 0)
    ,@('Tick', 4, 0)
    ,@('TickBitmap', 1, 0)
    ,@('TickMath', 0, 0)
    ,@('TransferHelper', 1, 0)
    ,@('UnsafeMath', 0, 0)
    ,@('NoDelegateCall', 0, 0)
    ,@('BitMathEchidnaTest', 1, 0)
    ,@('BitMathTest', 1, 0)
    ,@('FullMathEchidnaTest', 1, 0)
    ,@('FullMathTest', 1, 0)
    ,@('LiquidityMathTest', 1, 0)
    ,@('LowGasSafeMathEchidnaTest', 1, 0)
    ,@('MockTimeUniswapV3Pool', 1, 0)
    ,@('MockTimeUniswapV3PoolDeployer', 2, 0)
    ,@('NoDelegateCallTest', 1, 0)
    ,@('OracleEchidnaTest', 1, 0)
    ,@('OracleTest', 1, 0)
    ,@('SqrtPriceMathEchidnaTest', 3, 0)
    ,@('SqrtPriceMathTest', 1, 0)
    ,@('SwapMathEchidnaTest', 1, 0)
    ,@('SwapMathTest', 1, 0)
    ,@('TestERC20', 1, 0)
    ,@('TestUniswapV3Callee', 7, 0)
    ,@('TestUniswapV3ReentrantCallee', 3, 0)
    ,@('TestUniswapV3Router', 5, 0)
    ,@('TestUniswapV3SwapPay', 3, 0)
    ,@('TickBitmapEchidnaTest', 1, 0)
    ,@('TickBitmapTest', 1, 0)
    ,@('TickEchidnaTest', 1, 0)
    ,@('TickMathEchidnaTest', 1, 0)
    ,@('TickMathTest', 1, 0)
    ,@('TickOverflowSafetyEchidnaTest', 1, 0)
    ,@('TickTest', 1, 0)
    ,@('UniswapV3PoolSwapTest', 3, 0)
    ,@('UnsafeMathEchidnaTest', 1, 0)
    ,@('UniswapV3Factory', 4, 0)
    ,@('UniswapV3Pool', 21, 0)
    ,@('UniswapV3PoolDeployer', 2, 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Remove now-unused rows 64 to 77 (old data had 77 rows, new only has 63)
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -ge 64) {
    $deleteRange = $ws.Range("A64:C" + $lastRow)
    $deleteRange.EntireRow.Delete()
}

